$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number for each data row (rows 2-74).
# Update every row's value from 45177 to 45178 (2023-09-08 -> 2023-09-09),
# mirroring an automatic daily refresh of the "changed" timestamp.
for ($row = 2; $row -le 74; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45177) {
        $cell.Value = 45178
    }
}
